$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Hoja1" - the sheet with tabSelected="1"

# New WEB-related TODO items to append to rows 30-37 (column B),
# with column A style matching "not started" for rows 30-32.
$items = @(
    "WEB: Login",
    "WEB: Modificar datos (ya sea admin, profe o alumno)",
    "WEB: Estado academico (alumno)",
    "WEB: Inscripcion a materia (alumno)",
    "WEB: falta validar ""ModificarDatos.aspx""",
    "WEB: Cargar notas (profe)",
    "WEB: Inscribirse a curso (profe)",
    "WEB: ABMs del admin (son un monton!!)"
)

$startRow = 30
for ($i = 0; $i -lt $items.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $items[$i]
}

# Rows 30-32: column A style switches from the "red" status style to
# the "blue / not-started" status style used across most of the sheet
# (row 28's A cell already carries that style, and stays untouched by
# this edit, so copy its formatting across via a formats-only paste).
$ws.Cells.Item(28, 1).Copy()
for ($row = 30; $row -le 32; $row++) {
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats
}

# Scroll the sheet view down to A26 (was A16) while keeping the
# existing selection on B37.
$ws.Range("A26").Select()
$ws.Range("B37").Select()
